$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price (D) column data cells so that
# numeric-looking strings (e.g. "6.70", "8.60") keep their exact text
# representation instead of being normalized to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.876.97"
$ws.Range("E2").Value = "  +4.39%  "
$ws.Range("D3").Value = "2.697.21"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "581.99"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "148.84"
$ws.Range("E6").Value = "  +4.08%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "2.726.02"
$ws.Range("E9").Value = "  +4.52%  "
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "0.112"
$ws.Range("E11").Value = "  +6.57%  "
$ws.Range("D12").Value = "0.385"
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "3.199.56"
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "26.60"
$ws.Range("E15").Value = "  +8.83%  "
$ws.Range("D16").Value = "62.804.57"
$ws.Range("D17").Value = "0.0000149"
$ws.Range("E17").Value = "  +6.42%  "
$ws.Range("D18").Value = "2.722.28"
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("D19").Value = "11.90"
$ws.Range("E19").Value = "  +4.85%  "
$ws.Range("D20").Value = "4.86"
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("D21").Value = "361.41"
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").Value = "6.94"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "0.529"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "65.23"
$ws.Range("E25").Value = "  +2.50%  "
$ws.Range("D26").Value = "0.165"
$ws.Range("E26").Value = "  +3.54%  "
$ws.Range("D27").Value = "8.60"
$ws.Range("E27").Value = "  +7.42%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "2.01"
$ws.Range("E29").Value = "  +5.36%  "
$ws.Range("D30").Value = "0.0₃0849"
$ws.Range("E30").Value = "  +6.33%  "
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").Value = "  +10.38%  "
$ws.Range("D32").Value = "169.31"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.18"
$ws.Range("E34").Value = "  +19.81%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "20.48"
$ws.Range("E35").Value = "  +5.34%  "
$ws.Range("D36").Value = "4.72"
$ws.Range("E36").Value = "  +10.89%  "
$ws.Range("E37").Value = "  +7.28%  "
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").Value = "  +10.01%  "
$ws.Range("E39").Value = "  +19.16%  "
$ws.Range("D40").Value = "350.01"
$ws.Range("E40").Value = "  +11.99%  "
$ws.Range("D41").Value = "4.25"
$ws.Range("E41").Value = "  +9.23%  "
$ws.Range("D42").Value = "39.19"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("D43").Value = "5.59"
$ws.Range("E43").Value = "  +12.67%  "
$ws.Range("D44").Value = "21.52"
$ws.Range("E44").Value = "  +8.07%  "
$ws.Range("D45").Value = "0.0591"
$ws.Range("E45").Value = "  +7.47%  "
$ws.Range("D46").Value = "21.56"
$ws.Range("E46").Value = "  +8.44%  "
$ws.Range("D47").Value = "0.0259"
$ws.Range("E47").Value = "  +6.45%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "137.67"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.638"
$ws.Range("E49").Value = "  +5.29%  "
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "0.995"
$ws.Range("E51").Value = "  -0.39%  "
